$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5789666666666667
$ws.Range("H2").Value = 1.7369
$ws.Range("I2").Value = 0.01523705650035473
$ws.Range("J2").Value = 0.01523705650035472
$ws.Range("M2").Value = 211.980367
$ws.Range("N2").Value = 635.9411009999999
$ws.Range("O2").Value = 0.9885149156420702
$ws.Range("P2").Value = 0.9885149156420702
$ws.Range("Q2").Value = 122.7295664807667
$ws.Range("R2").Value = 1104.5660983269
$ws.Range("S2").Value = 0.01506205762108161
$ws.Range("T2").Value = 0.01506205762108161
$ws.Range("G3").Value = 0.5789666666666667
$ws.Range("H3").Value = 1.7369
$ws.Range("I3").Value = 0.01523705650035473
$ws.Range("J3").Value = 0.01523705650035472
$ws.Range("O3").Value = 0.003992992409159323
$ws.Range("P3").Value = 0.003992992409159324
$ws.Range("Q3").Value = 0.4957519806555556
$ws.Range("R3").Value = 4.461767825900001
$ws.Range("S3").Value = 0.00006084145094384814
$ws.Range("T3").Value = 0.00006084145094384814
$ws.Range("G4").Value = 0.5789666666666667
$ws.Range("H4").Value = 1.7369
$ws.Range("I4").Value = 0.01523705650035473
$ws.Range("J4").Value = 0.01523705650035472
$ws.Range("O4").Value = 0.007492091948770576
$ws.Range("P4").Value = 0.007492091948770576
$ws.Range("Q4").Value = 0.9301844437111112
$ws.Range("R4").Value = 8.371659993400002
$ws.Range("S4").Value = 0.00011415742832927
$ws.Range("T4").Value = 0.00011415742832927
$ws.Range("I5").Value = 0.6545086962501954
$ws.Range("J5").Value = 0.6545086962501954
$ws.Range("M5").Value = 211.980367
$ws.Range("N5").Value = 635.9411009999999
$ws.Range("O5").Value = 0.9885149156420702
$ws.Range("P5").Value = 0.9885149156420702
$ws.Range("Q5").Value = 5271.856053484359
$ws.Range("R5").Value = 47446.70448135924
$ws.Range("S5").Value = 0.6469916086607632
$ws.Range("T5").Value = 0.6469916086607632
$ws.Range("I6").Value = 0.6545086962501954
$ws.Range("J6").Value = 0.6545086962501954
$ws.Range("O6").Value = 0.003992992409159323
$ws.Range("P6").Value = 0.003992992409159324
$ws.Range("S6").Value = 0.002613448255855795
$ws.Range("T6").Value = 0.002613448255855795
$ws.Range("I7").Value = 0.6545086962501954
$ws.Range("J7").Value = 0.6545086962501954
$ws.Range("O7").Value = 0.007492091948770576
$ws.Range("P7").Value = 0.007492091948770576
$ws.Range("S7").Value = 0.004903639333576415
$ws.Range("T7").Value = 0.004903639333576415
$ws.Range("I8").Value = 0.33025424724945
$ws.Range("J8").Value = 0.3302542472494499
$ws.Range("M8").Value = 211.980367
$ws.Range("N8").Value = 635.9411009999999
$ws.Range("O8").Value = 0.9885149156420702
$ws.Range("P8").Value = 0.9885149156420702
$ws.Range("Q8").Value = 2660.091244815776
$ws.Range("R8").Value = 23940.82120334198
$ws.Range("S8").Value = 0.3264612493602254
$ws.Range("T8").Value = 0.3264612493602254
$ws.Range("I9").Value = 0.33025424724945
$ws.Range("J9").Value = 0.3302542472494499
$ws.Range("O9").Value = 0.003992992409159323
$ws.Range("P9").Value = 0.003992992409159324
$ws.Range("S9").Value = 0.00131870270235968
$ws.Range("T9").Value = 0.00131870270235968
$ws.Range("I10").Value = 0.33025424724945
$ws.Range("J10").Value = 0.3302542472494499
$ws.Range("O10").Value = 0.007492091948770576
$ws.Range("P10").Value = 0.007492091948770576
$ws.Range("S10").Value = 0.002474295186864891
$ws.Range("T10").Value = 0.002474295186864891
